$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 133, shifting existing rows 133-252 down to 134-253.
$ws.Rows(133).Insert()

# Populate the newly inserted row 133 with its data.
$ws.Range("A133").Value = 5
$ws.Range("B133").Value = "Macroferia Regional de Talca"
$ws.Range("C133").Value = "Maule"
$ws.Range("D133").Value = 44566
$ws.Range("E133").Value = 7
$ws.Range("F133").Value = 100112023
$ws.Range("G133").Value = "Brócoli"
$ws.Range("H133").Value = "Sin especificar"
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 4000
$ws.Range("K133").Value = 500
$ws.Range("L133").Value = 500
$ws.Range("M133").Value = 500
$ws.Range("N133").Value = "$/unidad"
$ws.Range("O133").Value = "Región del Maule"
$ws.Range("P133").Value = 500
$ws.Range("Q133").Value = 1
$ws.Range("R133").Value = "Hortaliza"
